$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# Copy the well-coding values ("NL_CW_W_11" .. "NL_CW_W_122") from the
# "Chemical analysis" sheet's column A (rows 3-24) into the matching rows
# of column A on the "In-situ measurements" sheet, replacing the previous
# "NL_CW_W_01".."NL_CW_W_022" coding values used there.
for ($r = 3; $r -le 24; $r++) {
    $v = $ws1.Cells.Item($r, 1).Value2
    $ws2.Cells.Item($r, 1).Value = $v
}

# Update the selection on the "Chemical analysis" sheet to A3:A24 (the
# range that was just copied from).
[void]$ws1.Activate()
[void]$ws1.Range("A3:A24").Select()

# Make the "In-situ measurements" sheet the active tab, with C26 selected.
[void]$ws2.Activate()
[void]$ws2.Range("C26").Select()
